$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value into a cell (prevents Excel from
# auto-converting numeric-looking strings like "580.90" into numbers),
# then restores the cell formatting so no stray number format lingers.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '67.765.30'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '3.260.08'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '580.90'
$ws.Range("E5").Value = '  -1.03%  '
Set-TextValue $ws.Range("D6") '185.26'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").Value = '3.260.96'
$ws.Range("E9").Value = '  -0.35%  '
Set-TextValue $ws.Range("D10") '0.131'
$ws.Range("E10").Value = '  -2.76%  '
Set-TextValue $ws.Range("D11") '6.58'
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").Value = '3.824.86'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  -0.25%  '
Set-TextValue $ws.Range("D15") '27.75'
$ws.Range("E15").Value = '  -3.02%  '
$ws.Range("D16").Value = '67.751.10'
$ws.Range("E16").Value = '  -0.61%  '
Set-TextValue $ws.Range("D17") '0.0000169'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '3.302.68'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("E19").Value = '  -1.98%  '
Set-TextValue $ws.Range("D20") '13.60'
$ws.Range("E20").Value = '  -0.08%  '
Set-TextValue $ws.Range("D21") '393.93'
$ws.Range("E21").Value = '  +3.00%  '
$ws.Range("E22").Value = '  -1.82%  '
Set-TextValue $ws.Range("D23") '71.56'
$ws.Range("E23").Value = '  +0.25%  '
Set-TextValue $ws.Range("D24") '0.999'
$ws.Range("E24").Value = '  -0.10%  '
Set-TextValue $ws.Range("D25") '0.517'
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("E27").Value = '  -2.62%  '
Set-TextValue $ws.Range("D28") '9.61'
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -1.84%  '
Set-TextValue $ws.Range("D31") '5.55'
$ws.Range("E31").Value = '  -4.41%  '
$ws.Range("E32").Value = '  -1.13%  '
Set-TextValue $ws.Range("D33") '6.99'
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("E34").Value = '  -3.00%  '
Set-TextValue $ws.Range("D36") '162.96'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  -3.73%  '
Set-TextValue $ws.Range("D38") '1.90'
$ws.Range("E38").Value = '  +1.67%  '
Set-TextValue $ws.Range("D39") '26.65'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  -3.44%  '
Set-TextValue $ws.Range("D41") '4.54'
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("E43").Value = '  -5.64%  '
$ws.Range("E44").Value = '  -0.06%  '
Set-TextValue $ws.Range("D45") '40.65'
$ws.Range("D46").Value = '2.614.98'
$ws.Range("E46").Value = '  -0.57%  '
Set-TextValue $ws.Range("D47") '24.89'
$ws.Range("E47").Value = '  -2.81%  '
Set-TextValue $ws.Range("D48") '335.12'
$ws.Range("E48").Value = '  -1.71%  '
$ws.Range("E49").Value = '  -2.05%  '
Set-TextValue $ws.Range("D50") '6.37'
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("E51").Value = '  -0.63%  '
